$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: a new price observation is inserted as row 152
# (sorted by date), pushing the existing rows 152-198 down to 153-199.
$ws.Rows.Item(152).Insert()

$ws.Range("A152").Value = 4
$ws.Range("B152").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C152").Value = "Los Lagos"
$ws.Range("D152").Value = 45093
$ws.Range("E152").Value = 10
$ws.Range("F152").Value = 100112052
$ws.Range("G152").Value = "Albahaca"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 90
$ws.Range("K152").Value = 5000
$ws.Range("L152").Value = 5000
$ws.Range("M152").Value = 5000
$ws.Range("N152").Value = "$/paquete"
$ws.Range("O152").Value = "Región de Arica y Parinacota"
$ws.Range("P152").Value = 5000
$ws.Range("Q152").Value = 1
$ws.Range("R152").Value = "Hortaliza"
